$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.360.25"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "'1.880.03"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'0.7097"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'242.54"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07998"
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("D9").Value = "'0.3136"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").Value = "'25.12"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "'0.08355"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").Value = "'1.878.39"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "'5.255"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'94.74"
$ws.Range("E14").Value = "  +3.51%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.7173"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "'6.311"
$ws.Range("E16").Value = "  +4.20%  "
$ws.Range("D17").Value = "'0.000008512"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").Value = "'29.379.51"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'242.30"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "'2.139.39"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "'13.26"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'7.828"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'0.1577"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").Value = "'163.60"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'9.086"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "'1.509"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").Value = "'4.423"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "'4.341"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  -6.29%  "
$ws.Range("D33").Value = "'0.05408"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "'1.939"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'0.7742"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("D36").Value = "'1.179"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'2.688"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'0.01883"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "'1.275.69"
$ws.Range("E39").Value = "  +6.29%  "
$ws.Range("D40").Value = "'2.745"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Value = "'6.579"
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("D42").Value = "'0.9264"
$ws.Range("E42").Value = "  +4.28%  "
$ws.Range("D43").Value = "'112.70"
$ws.Range("E43").Value = "  +5.05%  "
$ws.Range("D44").Value = "'74.63"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'0.00000000128"
$ws.Range("E46").Value = "  +4.22%  "
$ws.Range("D47").Value = "'2.034.34"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "'1.805"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").Value = "'0.5221"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").Value = "'9.550"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").Value = "'0.4366"
$ws.Range("E51").Value = "  +0.86%  "
